$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.734.47"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.601.91"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'211.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'19.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.826.35"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "1.599.33"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'65.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "'210.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'143.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "'7.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D31").Value = "'3.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "1.296.57"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "'2.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "'1.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").Value = "'0.601"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").Value = "'1.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.38%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'5.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'0.783"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'63.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "1.738.59"
$ws.Range("D45").Value = "'90.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'0.0516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'7.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "'0.397"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
